# Update the forecast rows (years 2022-2042, sheet rows 54-74) with the
# revised scenario values for columns B:H. Column A (year) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then values for columns B,C,D,E,F,G,H
$rows = @(
    @{ Row = 54; Vals = @(33.2822, 99.8466, 21.3957, 45.1687, 16.6411, 21.3957, 237.73) },
    @{ Row = 55; Vals = @(34.04799999999999, 102.144, 21.888, 46.208, 17.024, 21.888, 243.2) },
    @{ Row = 56; Vals = @(34.1628, 102.4884, 21.9618, 46.3638, 17.0814, 21.9618, 244.02) },
    @{ Row = 57; Vals = @(34.9244, 104.7732, 22.4514, 47.3974, 17.4622, 22.4514, 249.46) },
    @{ Row = 58; Vals = @(34.979, 104.937, 22.4865, 47.4715, 17.4895, 22.4865, 249.85) },
    @{ Row = 59; Vals = @(35.6062, 106.8186, 22.8897, 48.3227, 17.8031, 22.8897, 254.33) },
    @{ Row = 60; Vals = @(36.099, 108.297, 23.2065, 48.9915, 18.0495, 23.2065, 257.85) },
    @{ Row = 61; Vals = @(36.3902, 109.1706, 23.3937, 49.3867, 18.1951, 23.3937, 259.93) },
    @{ Row = 62; Vals = @(36.8886, 110.6658, 23.7141, 50.06310000000001, 18.4443, 23.7141, 263.49) },
    @{ Row = 63; Vals = @(37.31, 111.93, 23.985, 50.635, 18.655, 23.985, 266.5) },
    @{ Row = 64; Vals = @(37.6796, 113.0388, 24.2226, 51.1366, 18.8398, 24.2226, 269.14) },
    @{ Row = 65; Vals = @(38.2032, 114.6096, 24.5592, 51.8472, 19.1016, 24.5592, 272.88) },
    @{ Row = 66; Vals = @(38.5504, 115.6512, 24.7824, 52.3184, 19.2752, 24.7824, 275.36) },
    @{ Row = 67; Vals = @(38.9172, 116.7516, 25.0182, 52.81620000000001, 19.4586, 25.0182, 277.98) },
    @{ Row = 68; Vals = @(39.3456, 118.0368, 25.2936, 53.3976, 19.6728, 25.2936, 281.04) },
    @{ Row = 69; Vals = @(39.8174, 119.4522, 25.5969, 54.03790000000001, 19.9087, 25.5969, 284.41) },
    @{ Row = 70; Vals = @(40.1058, 120.3174, 25.78230000000001, 54.4293, 20.0529, 25.78230000000001, 286.47) },
    @{ Row = 71; Vals = @(40.5986, 121.7958, 26.0991, 55.0981, 20.2993, 26.0991, 289.99) },
    @{ Row = 72; Vals = @(40.8982, 122.6946, 26.2917, 55.5047, 20.4491, 26.2917, 292.13) },
    @{ Row = 73; Vals = @(41.3728, 124.1184, 26.5968, 56.14879999999999, 20.6864, 26.5968, 295.52) },
    @{ Row = 74; Vals = @(41.8012, 125.4036, 26.8722, 56.7302, 20.9006, 26.8722, 298.58) }
)

$cols = @("B", "C", "D", "E", "F", "G", "H")

foreach ($entry in $rows) {
    $r = $entry.Row
    $vals = $entry.Vals
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}
